# Sync attendance_reports: normalize "Recorded By" (column G) ordering
# so that any entry containing "System" has its comma-separated parts
# reversed (e.g. "user@x.com, System" -> "System, user@x.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    if ($trimmed -contains "System") {
        $reversedParts = @()
        for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
            $reversedParts += $trimmed[$i]
        }
        $newText = [string]::Join(", ", $reversedParts)
        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}
